$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 229 ("くつろぎなさい" entry), shifting all following rows up by one.
$ws.Rows.Item(229).Delete()
